$d = $word.ActiveDocument

# Replace first name: Сергей -> Иван
$d.Content.Find.Execute("Сергей", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Иван", 2)

# Replace patronymic: Юрьевич -> Иванович
$d.Content.Find.Execute("Юрьевич", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Иванович", 2)

# Replace surname: Ляшенко -> Иванов (preserve trailing space before tabs)
$d.Content.Find.Execute("Ляшенко ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Иванов ", 2)

# Replace phone line: мой номер 8524456. -> мой номер: 562564.
$d.Content.Find.Execute("мой номер 8524456.", $true, $false, $false, $false, $false,
                         $true, 1, $false, "мой номер: 562564.", 2)
